$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20 mirrors the structure/formatting of row 19 (date col A keeps
# the date number-format style; B:E are plain numbers).
$ws.Range("A19").Copy($ws.Range("A20"))
$ws.Range("A20").Value = 45986

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.622852459381209
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 2.447176337618551
